# Update build Manually. Add a colon in trap card data.
# Appends a full stop ("。") to the end of the "Remote bomb" (遥控炸弹)
# trap card's Chinese effect text in row 13, column D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D13").Value = "被弃置时：选弃牌堆1张怪物牌送墓。"

# Leave the selection where the author left it after the edit.
$ws.Range("D14").Select()
